# Applies the numeric corrections described by the commit diff to kev.constants.data.xlsx.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("equilibrium_concentrations")
$ws.Range("A2").Value = [double]"9.91922727314749e-08"
$ws.Range("B2").Value = [double]"0.000353446611903322"
$ws.Range("C2").Value = [double]"7.5337984131284e-07"
$ws.Range("D2").Value = [double]"8.25536517321998e-12"
$ws.Range("E2").Value = [double]"1.52588624774593e-07"
$ws.Range("A3").Value = [double]"0.000105333040224324"
$ws.Range("B3").Value = [double]"0.000108182220832239"
$ws.Range("C3").Value = [double]"0.00024486845629509"
$ws.Range("D3").Value = [double]"2.84932451856555e-06"
$ws.Range("E3").Value = [double]"1.43692923437208e-10"
$ws.Range("A4").Value = [double]"0.000358863895185699"
$ws.Range("B4").Value = [double]"3.9124914329389e-05"
$ws.Range("C4").Value = [double]"0.000301714024351818"
$ws.Range("D4").Value = [double]"1.19610613200632e-05"
$ws.Range("E4").Value = [double]"4.21764704875923e-11"
$ws.Range("A5").Value = [double]"0.000668387356554218"
$ws.Range("B5").Value = [double]"2.1694715995452e-05"
$ws.Range("C5").Value = [double]"0.000311597903157288"
$ws.Range("D5").Value = [double]"2.30073820908267e-05"
$ws.Range("E5").Value = [double]"2.26449712669488e-11"
$ws.Range("A6").Value = [double]"0.0013265284914767"
$ws.Range("B6").Value = [double]"1.05662358057948e-05"
$ws.Range("C6").Value = [double]"0.000301196008456199"
$ws.Range("D6").Value = [double]"4.41377557391445e-05"
$ws.Range("E6").Value = [double]"1.14099414988916e-11"
$ws.Range("A7").Value = [double]"0.00263114097144081"
$ws.Range("B7").Value = [double]"4.78545548604535e-06"
$ws.Range("C7").Value = [double]"0.000270570054716234"
$ws.Range("D7").Value = [double]"7.86444897977227e-05"
$ws.Range("E7").Value = [double]"5.75249013589486e-12"
$ws.Range("A8").Value = [double]"0.00661725897445892"
$ws.Range("B8").Value = [double]"1.43317099931708e-06"
$ws.Range("C8").Value = [double]"0.000203792630173007"
$ws.Range("D8").Value = [double]"0.000148974198827707"
$ws.Range("E8").Value = [double]"2.28729335556943e-12"
$ws.Range("A9").Value = [double]"0.0133758304235167"
$ws.Range("B9").Value = [double]"4.96251587698365e-07"
$ws.Range("C9").Value = [double]"0.000142637919209736"
$ws.Range("D9").Value = [double]"0.000210765829202615"
$ws.Range("E9").Value = [double]"1.13156432199914e-12"

$ws = $wb.Worksheets.Item("absorbance_calc_abs_errors")
$ws.Range("C2").Value = [double]"1.14941739155695"
$ws.Range("D2").Value = [double]"1.54289869907694"
$ws.Range("E2").Value = [double]"1.63781189285319"
$ws.Range("F2").Value = [double]"1.68384273214971"
$ws.Range("G2").Value = [double]"1.70338626130527"
$ws.Range("H2").Value = [double]"1.71146027647874"
$ws.Range("I2").Value = [double]"1.7434707888918"
$ws.Range("J2").Value = [double]"1.78869694814809"
$ws.Range("C3").Value = [double]"2.70361824233491"
$ws.Range("D3").Value = [double]"2.38169188018855"
$ws.Range("E3").Value = [double]"2.25461344910061"
$ws.Range("F3").Value = [double]"2.23501955066678"
$ws.Range("G3").Value = [double]"2.1835516640739"
$ws.Range("H3").Value = [double]"2.1115836822525"
$ws.Range("I3").Value = [double]"2.02114935259915"
$ws.Range("J3").Value = [double]"1.98878617183375"
$ws.Range("C4").Value = [double]"0.000417391556946756"
$ws.Range("D4").Value = [double]"-0.000101300923063707"
$ws.Range("E4").Value = [double]"-0.00418810714681084"
$ws.Range("F4").Value = [double]"-0.00115726785028913"
$ws.Range("G4").Value = [double]"0.00238626130526742"
$ws.Range("H4").Value = [double]"0.00746027647874192"
$ws.Range("I4").Value = [double]"-0.00652921110820381"
$ws.Range("J4").Value = [double]"0.00169694814809085"
$ws.Range("C5").Value = [double]"0.00161824233491226"
$ws.Range("D5").Value = [double]"-0.00730811981144797"
$ws.Range("E5").Value = [double]"0.00261344910061245"
$ws.Range("F5").Value = [double]"0.00701955066677984"
$ws.Range("G5").Value = [double]"-0.00244833592609517"
$ws.Range("H5").Value = [double]"-0.00241631774749651"
$ws.Range("I5").Value = [double]"0.00114935259914928"
$ws.Range("J5").Value = [double]"-0.000213828166249597"

$ws = $wb.Worksheets.Item("absorbance_calc_rel_errors")
$ws.Range("C2").Value = [double]"1.14941739155695"
$ws.Range("D2").Value = [double]"1.54289869907694"
$ws.Range("E2").Value = [double]"1.63781189285319"
$ws.Range("F2").Value = [double]"1.68384273214971"
$ws.Range("G2").Value = [double]"1.70338626130527"
$ws.Range("H2").Value = [double]"1.71146027647874"
$ws.Range("I2").Value = [double]"1.7434707888918"
$ws.Range("J2").Value = [double]"1.78869694814809"
$ws.Range("C3").Value = [double]"2.70361824233491"
$ws.Range("D3").Value = [double]"2.38169188018855"
$ws.Range("E3").Value = [double]"2.25461344910061"
$ws.Range("F3").Value = [double]"2.23501955066678"
$ws.Range("G3").Value = [double]"2.1835516640739"
$ws.Range("H3").Value = [double]"2.1115836822525"
$ws.Range("I3").Value = [double]"2.02114935259915"
$ws.Range("J3").Value = [double]"1.98878617183375"
$ws.Range("C4").Value = [double]"0.000363265062616846"
$ws.Range("D4").Value = [double]"-6.56519268073282e-05"
$ws.Range("E4").Value = [double]"-0.00255061336590185"
$ws.Range("F4").Value = [double]"-0.000686805845868922"
$ws.Range("G4").Value = [double]"0.00140285791021013"
$ws.Range("H4").Value = [double]"0.00437809652508329"
$ws.Range("I4").Value = [double]"-0.00373097777611647"
$ws.Range("J4").Value = [double]"0.000949607245713961"
$ws.Range("C5").Value = [double]"0.000598905379316157"
$ws.Range("D5").Value = [double]"-0.00305907066197068"
$ws.Range("E5").Value = [double]"0.00116050137682613"
$ws.Range("F5").Value = [double]"0.00315060622386887"
$ws.Range("G5").Value = [double]"-0.00112000728549642"
$ws.Range("H5").Value = [double]"-0.00114300744914688"
$ws.Range("I5").Value = [double]"0.000568986435222415"
$ws.Range("J5").Value = [double]"-0.000107505362619204"

$ws = $wb.Worksheets.Item("correlation_matrix")
$ws.Range("B2").Value = [double]"0.740919137027453"
$ws.Range("A3").Value = [double]"0.740919137027453"

$ws = $wb.Worksheets.Item("adj_r_squared")
$ws.Range("A2").Value = [double]"0.999842044335903"

$ws = $wb.Worksheets.Item("mol_ext_coefficients_calc")
$ws.Range("C2").Value = [double]"6.92174319352874"
$ws.Range("D2").Value = [double]"3241.77065205401"
$ws.Range("E2").Value = [double]"4810.07287563224"
$ws.Range("F2").Value = [double]"4784.48223182565"
$ws.Range("C3").Value = [double]"14.5560701290438"
$ws.Range("D3").Value = [double]"7635.87014621899"
$ws.Range("E3").Value = [double]"6297.39630669605"
$ws.Range("F3").Value = [double]"4232.42150409292"
$ws.Range("C4").Value = [double]"2.18811427811845"
$ws.Range("D4").Value = [double]"15.5429036014627"
$ws.Range("E4").Value = [double]"12.5548913967064"
$ws.Range("F4").Value = [double]"131.478092353871"
$ws.Range("C5").Value = [double]"2.18295624471658"
$ws.Range("D5").Value = [double]"15.5062643743712"
$ws.Range("E5").Value = [double]"12.5252957993336"
$ws.Range("F5").Value = [double]"131.168159550655"

$ws = $wb.Worksheets.Item("constants_evaluated")
$ws.Range("B4").Formula = '="4.33221233129501"'
$ws.Range("B4").Copy()
$ws.Range("B4").PasteSpecial(-4163)
$ws.Range("C4").Formula = '="0.0500521534968924"'
$ws.Range("C4").Copy()
$ws.Range("C4").PasteSpecial(-4163)
$ws.Range("B5").Formula = '="6.37545678377151"'
$ws.Range("B5").Copy()
$ws.Range("B5").PasteSpecial(-4163)
$ws.Range("C5").Formula = '="0.224859559792071"'
$ws.Range("C5").Copy()
$ws.Range("C5").PasteSpecial(-4163)
$excel.CutCopyMode = $false
